# Refresh the cryptocurrency price/volume table with the latest scrape.
# Price values that look numeric (e.g. "1.002", "306.22", "0.000008613")
# must be forced to Text so Excel stores them verbatim (matching the
# source data's inline-string cells, e.g. preserving trailing zeros such
# as "146.60") instead of silently coercing them to doubles. The
# NumberFormat is reset to the default "Normal" style right after the
# assignment so no stray text-format style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.172.55"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.902.81"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5261"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3776"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07259"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8992"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.36%  "
$ws.Range("D13").Value = "1.895.93"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.272"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008613"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "27.211.37"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.064"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "2.135.96"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.435"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.929"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.786"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09252"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8131"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05054"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.344"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.570"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5705"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01977"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.951"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4838"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.612"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
